# Append one new observation record as row 9 to the sheet (matches the
# diff: sheet dimension grows from A1:AY8 to A1:AY9, and a new row 9 is
# added with the data below).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9

# --- Numeric fields -------------------------------------------------------
$ws.Cells.Item($row, 1).Value  = 131153866   # A  Id
$ws.Cells.Item($row, 2).Value  = 58043       # B  Taxonsorteringsordning
$ws.Cells.Item($row, 5).Value  = 103021      # E  TaxonId
$ws.Cells.Item($row, 17).Value = 311201      # Q  Ost
$ws.Cells.Item($row, 18).Value = 6410515     # R  Nord
$ws.Cells.Item($row, 19).Value = 10          # S  Noggrannhet

# --- Plain text fields -----------------------------------------------------
$ws.Cells.Item($row, 4).Value  = "NT"                               # D  Rödlistade
$ws.Cells.Item($row, 6).Value  = "Talltita"                         # F  Artnamn
$ws.Cells.Item($row, 7).Value  = "Poecile montanus"                 # G  Vetenskapligt namn
$ws.Cells.Item($row, 8).Value  = "(Conrad von Baldenstein, 1827)"   # H  Auktor
$ws.Cells.Item($row, 16).Value = "Korseberget, Boh"                 # P  Lokalnamn
$ws.Cells.Item($row, 20).Value = "Västra Götaland"                  # T  Län
$ws.Cells.Item($row, 21).Value = "Kungälv"                          # U  Kommun
$ws.Cells.Item($row, 22).Value = "Bohuslän"                         # V  Provins
$ws.Cells.Item($row, 23).Value = "Harestad"                         # W  Socken
$ws.Cells.Item($row, 29).Value = "Inspelad i fält med en Audiomoth inspelningsapparat."  # AC Publik kommentar
$ws.Cells.Item($row, 49).Value = "Linus Lundin"                     # AW Rapportör
$ws.Cells.Item($row, 50).Value = "Linus Lundin"                     # AX Observatörer

# --- Date-like text fields (force text so Excel doesn't auto-convert them
#     into real date serials) --------------------------------------------
$ws.Cells.Item($row, 25).NumberFormat = "@"   # Y  Startdatum
$ws.Cells.Item($row, 25).Value = "2026-02-09"
$ws.Cells.Item($row, 25).Style = "Normal"     # drop the text number format again

$ws.Cells.Item($row, 27).NumberFormat = "@"   # AA Slutdatum
$ws.Cells.Item($row, 27).Value = "2026-02-09"
$ws.Cells.Item($row, 27).Style = "Normal"

# --- Boolean fields ---------------------------------------------------------
$ws.Cells.Item($row, 30).Value = $false   # AD Ej återfunnen
$ws.Cells.Item($row, 31).Value = $false   # AE Osäker artbestämning
$ws.Cells.Item($row, 33).Value = $false   # AG Ospontan

# --- Columns that are present in the source row but hold no value (typed,
#     empty text cells). A leading apostrophe forces an explicit empty text
#     value instead of leaving the cell absent; resetting the style back to
#     Normal drops the quote-prefix formatting that the apostrophe implies.
#     Columns: I, K, L, M, N, AT, AY
foreach ($c in 9, 11, 12, 13, 14, 46, 51) {
    $cell = $ws.Cells.Item($row, $c)
    $cell.Value = "'"
    $cell.Style = "Normal"
}
